$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range($ws.Cells.Item(15,2), $ws.Cells.Item(15,13)).Merge() | Out-Null

$b = $ws.Cells.Item(15,2)
$b.Borders.Item(7).LineStyle = 1
$b.Borders.Item(7).Color = 0
Write-Host "after border left set"
for ($col=2; $col -le 13; $col++) {
    Write-Host "col $col border-left LineStyle" $ws.Cells.Item(15,$col).Borders.Item(7).LineStyle
}
